# Generate Report for Handback
# Adds a "handback version mismatch" error result for the
# 29ea539e-33ee-4011-91ce-a54244851909 entry (row 7) on both the
# zh-cn and de-de worksheets: a new "Latest Target File" hyperlink,
# an updated "Latest Handback File" / "Latest Handback DateTime", a
# long "Error Detail" message, and widens the Error Detail column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2eee0cac68d655883aaf203c23fa19739228ee3a/e2e/29ea539e-33ee-4011-91ce-a54244851909.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc08c64a521f8bee6948381bebc98024a29d58ab/e2e/29ea539e-33ee-4011-91ce-a54244851909.md."
$latestMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc08c64a521f8bee6948381bebc98024a29d58ab/e2e/29ea539e-33ee-4011-91ce-a54244851909.md"
$mdDisplay = "29ea539e-33ee-4011-91ce-a54244851909.md"

# --- zh-cn sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Columns.Item(16).ColumnWidth = 39.17

$ws.Hyperlinks.Add($ws.Range("I7"), $latestMdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay) | Out-Null

$ws.Range("J7").Value = "29ea539e-33ee-4011-91ce-a54244851909.d238e3f52f410c87d19d12c611397162d0f7f7bd.zh-cn.xlf"
$ws.Range("K7").Value = "2016-08-28 02:42:45"
$ws.Range("P7").Value = $errorDetail

# --- de-de sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Columns.Item(16).ColumnWidth = 39.17

$ws.Hyperlinks.Add($ws.Range("I7"), $latestMdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay) | Out-Null

$ws.Range("J7").Value = "29ea539e-33ee-4011-91ce-a54244851909.d238e3f52f410c87d19d12c611397162d0f7f7bd.de-de.xlf"
$ws.Range("K7").Value = "2016-08-28 02:42:52"
$ws.Range("P7").Value = $errorDetail
